$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.450.27"
$ws.Range("E2").Value = "  -4.32%  "
$ws.Range("D3").Value = "2.367.37"
$ws.Range("E3").Value = "  -6.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.553"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("D9").Value = "2.386.26"
$ws.Range("E9").Value = "  -5.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0958"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("E12").Value = "  -9.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.319"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.81%  "
$ws.Range("D14").Value = "2.786.66"
$ws.Range("E14").Value = "  -5.21%  "
$ws.Range("D15").Value = "56.330.29"
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.71%  "
$ws.Range("D18").Value = "2.360.60"
$ws.Range("E18").Value = "  -6.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "312.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("E21").Value = "  -4.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("E26").Value = "  -4.90%  "
$ws.Range("D27").Value = "2.472.85"
$ws.Range("E27").Value = "  -5.45%  "
$ws.Range("E28").Value = "  -5.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.34%  "
$ws.Range("D33").Value = "0.0₃0713"
$ws.Range("E33").Value = "  -6.49%  "
$ws.Range("E34").Value = "  -7.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.13%  "
$ws.Range("E40").Value = "  -2.82%  "
$ws.Range("E41").Value = "  -5.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.787"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.64%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "129.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "254.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.571"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0902"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0488"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.53%  "
$ws.Range("E50").Value = "  -5.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.88%  "
